# Update countries & provincias Spain
# Refresh the COVID-19 stats table ("Pais" sheet) with the latest figures
# and bump the "datos actualizados" timestamp. Row positions / country
# names are unchanged - only the numeric columns (B:H) and the A1
# timestamp string move.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- timestamp banner (A1) ---------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 17 de Septiembre de 2020 a las 23:57"

# --- Estados Unidos (row 4) ---------------------------------------------
$ws.Range("B4").Value = 6867179
$ws.Range("C4").Value = 38878
$ws.Range("D4").Value = 4143865
$ws.Range("E4").Value = 2521208
$ws.Range("G4").Value = 758
$ws.Range("H4").Value = 202106

# --- Brasil (row 6) -------------------------------------------------------
$ws.Range("B6").Value = 4455386
$ws.Range("C6").Value = 33700
$ws.Range("E6").Value = 600139
$ws.Range("G6").Value = 761
$ws.Range("H6").Value = 134935

# --- Sudafrica (row 11) ---------------------------------------------------
$ws.Range("B11").Value = 655572
$ws.Range("C11").Value = 2128
$ws.Range("D11").Value = 585303
$ws.Range("E11").Value = 54497
$ws.Range("G11").Value = 67
$ws.Range("H11").Value = 15772

# --- Bulgaria (row 84) -----------------------------------------------------
$ws.Range("B84").Value = 18544
$ws.Range("C84").Value = 154
$ws.Range("D84").Value = 13391
$ws.Range("E84").Value = 4404
$ws.Range("G84").Value = 10
$ws.Range("H84").Value = 749

# --- Haiti (row 104) --------------------------------------------------------
$ws.Range("B104").Value = 8556
$ws.Range("C104").Value = 15
$ws.Range("D104").Value = 6315
$ws.Range("E104").Value = 2021

# --- Guinea Ecuatorial (row 116) --------------------------------------------
$ws.Range("B116").Value = 5063
$ws.Range("C116").Value = 85
$ws.Range("D116").Value = 4465
$ws.Range("E116").Value = 549
$ws.Range("G116").Value = 2
$ws.Range("H116").Value = 49

# --- Hong Kong (row 117) -----------------------------------------------------
$ws.Range("B117").Value = 5002
$ws.Range("C117").Value = 2
$ws.Range("D117").Value = 4509
$ws.Range("E117").Value = 410
$ws.Range("H117").Value = 83

# --- Cabo Verde (row 118) ----------------------------------------------------
$ws.Range("B118").Value = 4994
$ws.Range("C118").Value = 9
$ws.Range("D118").Value = 4682
$ws.Range("E118").Value = 210
$ws.Range("H118").Value = 102

# --- Yemen (row 152) ----------------------------------------------------------
$ws.Range("B152").Value = 2022
$ws.Range("C152").Value = 3
$ws.Range("E152").Value = 216
$ws.Range("G152").Value = 2
$ws.Range("H152").Value = 585

# --- Burkina Faso (row 156) ----------------------------------------------------
$ws.Range("B156").Value = 1767
$ws.Range("C156").Value = 19
$ws.Range("D156").Value = 1166
$ws.Range("E156").Value = 545

# --- Togo (row 157) --------------------------------------------------------------
$ws.Range("B157").Value = 1618
$ws.Range("C157").Value = 10
$ws.Range("D157").Value = 1243
$ws.Range("E157").Value = 334
$ws.Range("G157").Value = 1
$ws.Range("H157").Value = 41

# --- Republica de Chipre (row 158) ------------------------------------------------
$ws.Range("B158").Value = 1558
$ws.Range("C158").Value = 10
$ws.Range("E158").Value = 254

# --- Polinesia Francesa (row 166) -------------------------------------------------
$ws.Range("B166").Value = 1115
$ws.Range("C166").Value = 25
$ws.Range("D166").Value = 962
$ws.Range("E166").Value = 72
$ws.Range("H166").Value = 81

# --- Republica del Chad (row 167) -------------------------------------------------
$ws.Range("B167").Value = 1099
$ws.Range("D167").Value = 672
$ws.Range("E167").Value = 425
$ws.Range("H167").Value = 2

# --- Montserrat (row 214) -----------------------------------------------------------
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0

# --- Islas Malvinas (row 215) --------------------------------------------------------
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1
